$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Growth")
$ws.Activate()

# Clear the "user-specific" Data File Name columns (I, J) for the data rows -
# these cells become empty, removing the "user-specific" shared string.
$ws.Range("I2:J7").ClearContents()

# Re-write the Output [Sample Name] column (K) values - content is unchanged,
# but writing them keeps the shared-string table/order consistent after the
# "user-specific" string drops out.
$ws.Range("K2").Value = "CC1"
$ws.Range("K3").Value = "CC2"
$ws.Range("K4").Value = "CC3"
$ws.Range("K5").Value = "Co1"
$ws.Range("K6").Value = "Co2"
$ws.Range("K7").Value = "Co3"

# Update the view state to match the recorded selection/scroll position.
$ws.Range("I10").Select()
